# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The sheet's column G (header "K") holds per-game strikeout-type counts that were
# recomputed from source data; write the newly computed values for rows 2..64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(2,1,3,4,0,0,1,3,0,1,0,1,1,1,0,0,0,1,1,2,1,2,0,0,0,1,1,0,1,1,2,2,1,2,0,1,1,0,0,0,1,1,1,2,2,3,1,0,2,1,2,2,2,0,2,0,3,3,1,0,2,1,0)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
